$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "extra long basketball shorts for men"
$ws.Cells.Item(2, 1).Value = "compression pants running"
$ws.Cells.Item(3, 1).Value = "boys baseball compression sleeve"
$ws.Cells.Item(4, 1).Value = "stretch mark for men"
$ws.Cells.Item(5, 1).Value = "compression shorts men pack"
$ws.Cells.Item(6, 1).Value = "knee pads youth girls"
$ws.Cells.Item(7, 1).Value = "knee sleeves basketball"
$ws.Cells.Item(8, 1).Value = "xl youth baseball pants"
$ws.Cells.Item(9, 1).Value = "biking knee sleeve"
$ws.Cells.Item(10, 1).Value = "boy baseball pants"
$ws.Cells.Item(11, 1).Value = "men black compression pants"
$ws.Cells.Item(12, 1).Value = "knee pads exercise"
$ws.Cells.Item(13, 1).Value = "youth basketball compression shorts"
$ws.Cells.Item(14, 1).Value = "hockey pads"
$ws.Cells.Item(15, 1).Value = "calf strain compression sleeve"
$ws.Cells.Item(16, 1).Value = "mens degree sport"
$ws.Cells.Item(17, 1).Value = "men long shorts below the knee"
$ws.Cells.Item(18, 1).Value = "knee pads construction gel"
$ws.Cells.Item(19, 1).Value = "work knee pads construction"
$ws.Cells.Item(20, 1).Value = "padded football sleeve"
$ws.Cells.Item(21, 1).Value = "spandex tights"
$ws.Cells.Item(22, 1).Value = "adult black baseball pants"
$ws.Cells.Item(23, 1).Value = "girls baseball pants"
$ws.Cells.Item(24, 1).Value = "mens small black baseball pants"
$ws.Cells.Item(25, 1).Value = "youth softball compression sleeve"
$ws.Cells.Item(26, 1).Value = "padded shorts football"
$ws.Cells.Item(27, 1).Value = "3/4 shorts for men"
$ws.Cells.Item(28, 1).Value = "silicone strips for stretch marks"
$ws.Cells.Item(29, 1).Value = "calf compression sleeve basketball"
$ws.Cells.Item(30, 1).Value = "gel knee pads for construction"
$ws.Cells.Item(31, 1).Value = "youth baseball short pants"
$ws.Cells.Item(32, 1).Value = "exercise knee pad"
$ws.Cells.Item(33, 1).Value = "girls basketball shorts size 6"
$ws.Cells.Item(34, 1).Value = "mens short tights"
$ws.Cells.Item(35, 1).Value = "boys knee length shorts"
$ws.Cells.Item(36, 1).Value = "rodilleras de volleyball"
$ws.Cells.Item(37, 1).Value = "rodilleras volleyball"
$ws.Cells.Item(38, 1).Value = "compression running pants men"
$ws.Cells.Item(39, 1).Value = "knee pads cheap"
$ws.Cells.Item(40, 1).Value = "exercise knee pads"
$ws.Cells.Item(41, 1).Value = "pad knee"
$ws.Cells.Item(42, 1).Value = "compression workout pants men"
$ws.Cells.Item(43, 1).Value = "compression sleeve for bursitis"
$ws.Cells.Item(44, 1).Value = "mens running tights shorts"
$ws.Cells.Item(45, 1).Value = "mens compression girdle"
$ws.Cells.Item(46, 1).Value = "knee pads outdoor"
$ws.Cells.Item(47, 1).Value = "motorcycle knee pads for men"
$ws.Cells.Item(48, 1).Value = "water pants men"
$ws.Cells.Item(49, 1).Value = "pant baseball men"
$ws.Cells.Item(50, 1).Value = "volleyball knee pad"
$ws.Cells.Item(51, 1).Value = "knee pads for volleyball"
$ws.Cells.Item(52, 1).Value = "pro tights men"
$ws.Cells.Item(53, 1).Value = "sleeve knee pads"
$ws.Cells.Item(54, 1).Value = "womens compression leggings"
$ws.Cells.Item(55, 1).Value = "airsoft knee pads"
$ws.Cells.Item(56, 1).Value = "skateboard knee pads"
$ws.Cells.Item(57, 1).Value = "starter youth compression pants"
$ws.Cells.Item(58, 1).Value = "trolls knee pads"
$ws.Cells.Item(59, 1).Value = "mens compression pants marvel"
$ws.Cells.Item(60, 1).Value = "mouthguard basketball youth"
$ws.Cells.Item(61, 1).Value = "mcdavid youth knee pads"
$ws.Cells.Item(62, 1).Value = "nike compression leggings"
$ws.Cells.Item(63, 1).Value = "jordan flight mens basketball pants"
$ws.Cells.Item(64, 1).Value = "athletic compression pants"
$ws.Cells.Item(65, 1).Value = "nike pro compression leggings men"
$ws.Cells.Item(66, 1).Value = "emoji knee pads"
$ws.Cells.Item(67, 1).Value = "mens compression tights 3 4"
$ws.Cells.Item(68, 1).Value = "mens compression tights nike"
$ws.Cells.Item(69, 1).Value = "mens compression tights white"
$ws.Cells.Item(70, 1).Value = "eastbay compression pants"
$ws.Cells.Item(71, 1).Value = "women compression leggings"
$ws.Cells.Item(72, 1).Value = "ladies compression pants"
$ws.Cells.Item(73, 1).Value = "womans compression leggings"
$ws.Cells.Item(74, 1).Value = "mcdavid compression pants"
$ws.Cells.Item(75, 1).Value = "knee brace basketball youth"
$ws.Cells.Item(76, 1).Value = "basketball kids knee pads"
$ws.Cells.Item(77, 1).Value = "basketball knee pads mcdavid"
$ws.Cells.Item(78, 1).Value = "basketball youth jersey"
$ws.Cells.Item(79, 1).Value = "morris compression knee pads"
$ws.Cells.Item(80, 1).Value = "wonens compression leggings"
$ws.Cells.Item(81, 1).Value = "youth baketball knee pads"
$ws.Cells.Item(82, 1).Value = "olympic mens basketball"
$ws.Cells.Item(83, 1).Value = "track leggings men"
$ws.Cells.Item(84, 1).Value = "nike pro dry mens basketball tights"
$ws.Cells.Item(85, 1).Value = "goalie knee protectors"
$ws.Cells.Item(86, 1).Value = "basketball clothes for men"
$ws.Cells.Item(87, 1).Value = "mens running thermal compression pants"
$ws.Cells.Item(88, 1).Value = "capri pants for men adidas"
$ws.Cells.Item(89, 1).Value = "men running tights"
$ws.Cells.Item(90, 1).Value = "mens running tight"
$ws.Cells.Item(91, 1).Value = "men tights nike"
$ws.Cells.Item(92, 1).Value = "men tights short"
$ws.Cells.Item(93, 1).Value = "mens tights dance"
$ws.Cells.Item(94, 1).Value = "men running tight"
$ws.Cells.Item(95, 1).Value = "men tights green"
$ws.Cells.Item(96, 1).Value = "men tights pack"
$ws.Cells.Item(97, 1).Value = "mens tights grey"
$ws.Cells.Item(98, 1).Value = "mens tights nike"
$ws.Cells.Item(99, 1).Value = "mens tights pink"
$ws.Cells.Item(100, 1).Value = "men legging nike"
